$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price-table refresh.
# Columns B/C/D/E hold literal text (inline strings) in the source data,
# including price strings like "1.00" or "34.453.48" that Excel would
# otherwise auto-convert to numbers. Force Text format on Price (D) cells
# before writing so the stored content matches exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.453.48"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.803.06"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.67"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.33"
$ws.Range("E8").Value = "  +3.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.062.95"
$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.52"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.799.12"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("E16").Value = "  +3.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.420.99"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.09"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.06"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("E20").Value = "  -1.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.56"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.16"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.15"
$ws.Range("E25").Value = "  +3.68%  "

$ws.Range("E26").Value = "  +7.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.89"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("E28").Value = "  +1.08%  "

$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.84"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("E33").Value = "  -0.63%  "

$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.393.32"
$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.673"
$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("E37").Value = "  -6.57%  "

$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0189"
$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("E41").Value = "  -5.33%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.958"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +0.67%  "

$ws.Range("E44").Value = "  +8.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.38"
$ws.Range("E45").Value = "  -4.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.04"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.964.43"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.10"
$ws.Range("E50").Value = "  -1.76%  "

$ws.Range("E51").Value = "  -1.95%  "
